# Applies the "Experimental Strategy and Study Data types - 13 Test cases" edit:
# - Rewrites the 4 Cypher queries (stat/CALL query + 3 per-tab queries) to the
#   new apoc.coll.flatten/CALL-subquery based versions.
# - Bumps the workbook's query font from 12pt to 14pt (applied via the shared
#   "Normal" style so every cellXfs record that referenced the old font now
#   reflects the new size, mirroring how Excel reuses xf/font records).
# - Rows 2-4 grow (text wrap + bigger font) up to Excel's 409.5pt row-height
#   cap; the view's active selection moves from B5 to D4 with row 2 scrolled
#   to the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the four Cypher query cells --------------------------------

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "WGA" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "WGA" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "WGA" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "WGA" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "WGA" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id LIMIT 100
'@

$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "WGA" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name LIMIT 100
'@

# Column C on every data row shares the same stat/CALL query.
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Column B holds the per-tab query (Participants / Samples / Files).
$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# --- 2. Bump the query font from 12pt to 14pt ------------------------------
# The workbook's shared "Normal" style carries the 12pt font used by every
# formatted cell on the sheet; resizing it cascades to all the cellXfs that
# reference it (header cells, the query cells, and the trailing blank cells),
# exactly like selecting the whole sheet and bumping the font size in the
# Home ribbon.
$wb.Styles.Item("Normal").Font.Size = 14

# --- 3. Row heights grow to fit the longer, bigger-font query text --------
# Excel's row height cap is 409.5 points; the new query text wraps onto
# enough lines in columns B/C that every data row now hits that ceiling.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# --- 4. View state: selection moves to D4, row 2 scrolled to the top ------
$ws.Range("D4").Select()
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older hosts may not expose ActiveWindow scrolling; selection already
    # moved above, so this is a best-effort nicety only.
}
